# I-O Distribution Key.xlsx -- fill in the "Goes to:" (column C) destination
# line numbers on the "Hungarian P&L" sheet, shrink columns A-C to fit the
# new narrower layout, and switch the active/selected sheet from
# "Hungarian BS" to "Hungarian P&L".
$wb   = $excel.ActiveWorkbook
$wsBS = $wb.Worksheets.Item("Hungarian BS")
$wsPL = $wb.Worksheets.Item("Hungarian P&L")

# Column C ("Goes to" line numbers) for the P&L sheet
$wsPL.Range("C3").Value  = 10
$wsPL.Range("C4").Value  = 12
$wsPL.Range("C6").Value  = 16
$wsPL.Range("C7").Value  = 16
$wsPL.Range("C9").Value  = 18
$wsPL.Range("C11").Value = 24
$wsPL.Range("C12").Value = 26
$wsPL.Range("C13").Value = 32
$wsPL.Range("C14").Value = 24
$wsPL.Range("C15").Value = 26
$wsPL.Range("C17").Value = 30
$wsPL.Range("C18").Value = 30
$wsPL.Range("C19").Value = 30
$wsPL.Range("C21").Value = 42
$wsPL.Range("C22").Value = 32
$wsPL.Range("C23").Value = 44
$wsPL.Range("C25").Value = 52
$wsPL.Range("C26").Value = 54
$wsPL.Range("C27").Value = 52
$wsPL.Range("C28").Value = 54
$wsPL.Range("C29").Value = 56
$wsPL.Range("C30").Value = 58
$wsPL.Range("C31").Value = 60
$wsPL.Range("C32").Value = 62
$wsPL.Range("C33").Value = 66
$wsPL.Range("C36").Value = 70
$wsPL.Range("C38").Value = 76
$wsPL.Range("C40").Value = 72
$wsPL.Range("C41").Value = 74
$wsPL.Range("C42").Value = 78
$wsPL.Range("C43").Value = 78
$wsPL.Range("C48").Value = 92

# Narrower columns now that column C carries short numeric values
$wsPL.Columns.Item(1).ColumnWidth = 36.67
$wsPL.Columns.Item(2).ColumnWidth = 33.83
$wsPL.Columns.Item(3).ColumnWidth = 8.5

# Make "Hungarian P&L" the active/selected sheet (was "Hungarian BS")
$wsPL.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$wsPL.Range("C50").Select() | Out-Null
